$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels and value to include extra clarifying text
$ws.Range("C3").Value = "1"
$ws.Range("B1").Value = "年级（年份）"
$ws.Range("C1").Value = "班级(班级号)"
$ws.Range("A1").Value = "专业（专业代码）"

# Widen column A to fit the new, longer header text
$ws.Columns.Item(1).ColumnWidth = 27.285714285714285

# Restore the selection to D22 as in the saved workbook
$ws.Range("D22").Select()
